# Update out of stock parts
# The LCSC part numbers for the BLU (D1) and GRN (D2) LEDs were changed
# because the originally specified parts went out of stock.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("F405_PILL_BOM")

# D1 (BLU LED) - LCSC part number C72041 -> C2689219
$ws.Range("D12").Value = "C2689219"

# D2 (GRN LED) - LCSC part number C72043 -> C87326
$ws.Range("D13").Value = "C87326"
